$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 6 - this shifts the existing rows 6 and 7
# (Carolina Castro / Secretaria, Jose Gonzalez / Piloto) down to rows 7 and 8.
$ws.Rows.Item(6).Insert()

# Row 8 now holds the "Jose Gonzalez / Piloto" record (shifted down from the
# old row 7). Copy it into the newly blank row 6 so the text-typed cells keep
# their shared-string typing (no accidental date auto-conversion on column H).
$ws.Range("A8:K8").Copy()
$ws.Range("A6:K6").PasteSpecial()

# Fix up the sequential id columns for the new row 6 (it is user id 1).
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 1

$excel.CutCopyMode = $false

# Match the saved selection from the edited workbook.
$ws.Range("C7").Select() | Out-Null
